# Saldo.xlsx update:
#  - account 004384167 (DOUGLAS) balance: 9248.16 -> 1874.01
#  - account 004361159 (HFR)     balance: 438.68  -> 103401.94
#  - the sheet is kept sorted by the "Saldo" column (C), descending,
#    so after the balance changes the data rows are re-sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the two account rows by their account number (column A) so the
# script doesn't depend on the current (pre-sort) row ordering.
$douglasCell = $ws.Columns.Item(1).Find("004384167")
$hfrCell = $ws.Columns.Item(1).Find("004361159")

$douglasRow = $douglasCell.Row
$hfrRow = $hfrCell.Row

$ws.Cells.Item($douglasRow, 3).Value = 1874.01
$ws.Cells.Item($hfrRow, 3).Value = 103401.94

# The data body runs from row 2 (first account row, right under the
# "Conta"/"Nome"/"Saldo" header) down to the last contiguous row before
# the blank separator row that precedes the "Filtros aplicados" footer.
$headerRow = 1
$firstDataRow = $headerRow + 1
$lastDataRow = $ws.Cells.Item($firstDataRow, 1).End(4).Row

$sortRange = $ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($lastDataRow, 3))
$sortKey = $ws.Range($ws.Cells.Item($firstDataRow, 3), $ws.Cells.Item($lastDataRow, 3))

# xlDescending = 2
$sortRange.Sort($sortKey, 2)

Write-Host "Updated DOUGLAS (row $douglasRow) and HFR (row $hfrRow); re-sorted rows $firstDataRow..$lastDataRow by Saldo desc."
